$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook used to hard-code exactly 3 "Placa_*" blocks of 5 rows each
# (16 rows total incl. header). The app now accepts an arbitrary number of
# .res files, so the per-plate header rows ("Placa_0 1", "Placa_1 1",
# "Placa_3 1") are grouped at the top and the numbered sample rows are
# renumbered 2..7 in a single contiguous block (12 rows total incl. header).
# The Re (column D) values were also recomputed (x10, unit-scale fix).

# row -> (A label, B, C, D, E, F); labels that look like plain integers must
# be forced back to Text, otherwise Excel's Range.Value setter infers them
# as numbers.
$rows = @(
    @(2,  "Placa_0 1", 0.5171329573658486, 2.996145142857143, 40210.8610295068,  29.48833428571428, 0.09439878676581151),
    @(3,  "2",         0.5709444553684897, 3.308152777777777, 44025.51001771774, 36.4758275,         0.1039381158842682),
    @(4,  "3",         0.7585202529229677, 4.395085135135135, 58355.83396856042, 64.48898486486488,  0.1992425198041262),
    @(5,  "Placa_1 1", 0.5171329573658486, 2.996145142857143, 40210.8610295068,  29.48833428571428,  0.09439878676581151),
    @(6,  "Placa_3 1", 0.5171329573658486, 2.996145142857143, 40210.8610295068,  29.48833428571428,  0.09439878676581151),
    @(7,  "2",         0.5709444553684897, 3.308152777777777, 44025.51001771774, 36.4758275,         0.1039381158842682),
    @(8,  "3",         0.7585202529229677, 4.395085135135135, 58355.83396856042, 64.48898486486488,  0.1992425198041262),
    @(9,  "4",         0.9473532113458141, 5.489527222222224, 72426.23396458979, 100.8328161111111,  0.2417405763006262),
    @(10, "5",         1.140045356010704,  6.606468285714286, 86571.73739163704, 145.9750517142857,  0.3444630314620629),
    @(11, "6",         1.140045356010704,  6.606468285714286, 86571.73739163704, 145.9750517142857,  0.3444630314620629),
    @(12, "7",         1.140045356010704,  6.606468285714286, 86571.73739163704, 145.9750517142857,  0.3444630314620629)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $label  = $r[1]
    $aCell  = $ws.Cells.Item($rowNum, 1)

    # Plate-header labels ("Placa_n n") are already non-numeric text and
    # don't need help; purely-numeric labels need the column forced to
    # Text first so "2" etc. isn't silently re-typed as the number 2.
    if ($label -match '^-?[0-9]+(\.[0-9]+)?$') {
        $aCell.NumberFormat = "@"
    }
    $aCell.Value = $label

    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# The old layout had 16 data rows (3 plates x 5 samples + header); the new
# one only needs 12 (3 plate headers + 7 renumbered sample rows + header).
# Drop the now-unused tail rows.
$ws.Range("A13:F16").Delete()
